$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 9961.25
$ws.Range("I40").Value = 8555.799999999999
$ws.Range("J40").Value = 10600.091
$ws.Range("K40").Value = 8555.799999999999
$ws.Range("L40").Value = 10600.091
$ws.Range("M40").Value = -8380.799999999999
$ws.Range("N40").Value = -10950.091
$ws.Range("H70").Value = 1893.6666
$ws.Range("I70").Value = 957.4
$ws.Range("K70").Value = 2872.2
$ws.Range("M70").Value = -2602.2
$ws.Range("H73").Value = 1893.6666
$ws.Range("I73").Value = 957.4
$ws.Range("K73").Value = 2872.2
$ws.Range("M73").Value = -1936.2
$ws.Range("H132").Value = 4057.077
$ws.Range("I132").Value = 2728.5
$ws.Range("K132").Value = 8185.5
$ws.Range("M132").Value = -5655.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8115.2666
$ws.Range("J61").Value = 14799.5
$ws.Range("L61").Value = 14799.5
$ws.Range("N61").Value = -15223.5
$ws.Range("H74").Value = 8775155
$ws.Range("I74").Value = 11496771
$ws.Range("J74").Value = 5500.778
$ws.Range("K74").Value = 11496771
$ws.Range("L74").Value = 5500.778
$ws.Range("M74").Value = -11495897
$ws.Range("N74").Value = -7248.778
$ws.Range("H77").Value = 8775155
$ws.Range("I77").Value = 11496771
$ws.Range("J77").Value = 5500.778
$ws.Range("K77").Value = 57483855
$ws.Range("L77").Value = 27503.89
$ws.Range("M77").Value = -57479487
$ws.Range("N77").Value = -36239.89
$ws.Range("H88").Value = 4123.25
$ws.Range("J88").Value = 2647.9
$ws.Range("L88").Value = 2647.9
$ws.Range("N88").Value = -3459.9
$ws.Range("H91").Value = 4123.25
$ws.Range("J91").Value = 2647.9
$ws.Range("L91").Value = 2647.9
$ws.Range("N91").Value = -5455.9
$ws.Range("H97").Value = 3220
$ws.Range("I97").Value = 2240.4092
$ws.Range("J97").Value = 7530.2
$ws.Range("K97").Value = 2240.4092
$ws.Range("L97").Value = 7530.2
$ws.Range("M97").Value = -1744.4092
$ws.Range("N97").Value = -8522.200000000001
$ws.Range("H132").Value = 6683.148
$ws.Range("I132").Value = 3795.889
$ws.Range("J132").Value = 12457.667
$ws.Range("K132").Value = 11387.667
$ws.Range("L132").Value = 37373.001
$ws.Range("M132").Value = -8857.667000000001
$ws.Range("N132").Value = -42433.001
$ws.Range("H136").Value = 8115.2666
$ws.Range("J136").Value = 14799.5
$ws.Range("L136").Value = 44398.5
$ws.Range("N136").Value = -49498.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6199.2856
$ws.Range("I86").Value = 2097.25
$ws.Range("K86").Value = 2097.25
$ws.Range("M86").Value = -974.25
$ws.Range("H89").Value = 6199.2856
$ws.Range("I89").Value = 2097.25
$ws.Range("K89").Value = 10486.25
$ws.Range("M89").Value = -4870.25
$ws.Range("H94").Value = 976.4286
$ws.Range("I94").Value = 976.4286
$ws.Range("K94").Value = 976.4286
$ws.Range("M94").Value = -525.4286
$ws.Range("H105").Value = 37124.145
$ws.Range("J105").Value = 12979.75
$ws.Range("L105").Value = 12979.75
$ws.Range("N105").Value = -16473.75
$ws.Range("H107").Value = 2418.9546
$ws.Range("I107").Value = 2064.0527
$ws.Range("K107").Value = 2064.0527
$ws.Range("M107").Value = -144.0527000000002
$ws.Range("H134").Value = 3487.5833
$ws.Range("I134").Value = 1667.125
$ws.Range("K134").Value = 5001.375
$ws.Range("M134").Value = -2466.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26840.744
$ws.Range("I31").Value = 2749.1667
$ws.Range("K31").Value = 2749.1667
$ws.Range("M31").Value = -2454.1667
$ws.Range("H34").Value = 26840.744
$ws.Range("I34").Value = 2749.1667
$ws.Range("K34").Value = 2749.1667
$ws.Range("M34").Value = -2547.1667
$ws.Range("H107").Value = 988.3333
$ws.Range("I107").Value = 529.25
$ws.Range("J107").Value = 2824.6667
$ws.Range("K107").Value = 529.25
$ws.Range("L107").Value = 2824.6667
$ws.Range("M107").Value = 1390.75
$ws.Range("N107").Value = -6664.6667
$ws.Range("H132").Value = 4224.24
$ws.Range("J132").Value = 6365.8887
$ws.Range("L132").Value = 19097.6661
$ws.Range("N132").Value = -24157.6661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 197.4
$ws.Range("I10").Value = 53.42857
$ws.Range("J10").Value = 533.3333
$ws.Range("K10").Value = 160.28571
$ws.Range("L10").Value = 1599.9999
$ws.Range("M10").Value = -21.28570999999999
$ws.Range("N10").Value = -1877.9999
$ws.Range("H54").Value = 8978.200000000001
$ws.Range("I54").Value = 6999
$ws.Range("K54").Value = 20997
$ws.Range("M54").Value = -20438
$ws.Range("H75").Value = 5583.0835
$ws.Range("J75").Value = 5697.1
$ws.Range("L75").Value = 17091.3
$ws.Range("N75").Value = -19087.3
$ws.Range("H78").Value = 5583.0835
$ws.Range("J78").Value = 5697.1
$ws.Range("L78").Value = 51273.9
$ws.Range("N78").Value = -61257.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1238.4
$ws.Range("I97").Value = 855.5
$ws.Range("J97").Value = 2131.8333
$ws.Range("K97").Value = 855.5
$ws.Range("L97").Value = 2131.8333
$ws.Range("M97").Value = -359.5
$ws.Range("N97").Value = -3123.8333
$ws.Range("H102").Value = 2778.0476
$ws.Range("I102").Value = 1732.8334
$ws.Range("K102").Value = 1732.8334
$ws.Range("M102").Value = -110.8334
$ws.Range("H113").Value = 5065.4443
$ws.Range("I113").Value = 1264.3334
$ws.Range("J113").Value = 8866.556
$ws.Range("K113").Value = 1264.3334
$ws.Range("L113").Value = 8866.556
$ws.Range("M113").Value = 905.6666
$ws.Range("N113").Value = -13206.556
$ws.Range("H132").Value = 5466.0713
$ws.Range("I132").Value = 3783.5
$ws.Range("J132").Value = 8494.700000000001
$ws.Range("K132").Value = 11350.5
$ws.Range("L132").Value = 25484.1
$ws.Range("M132").Value = -8820.5
$ws.Range("N132").Value = -30544.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 19500
$ws.Range("I20").Value = 19500
$ws.Range("K20").Value = 19500
$ws.Range("M20").Value = -19274
$ws.Range("H22").Value = 33869
$ws.Range("I22").Value = 51636.5
$ws.Range("K22").Value = 51636.5
$ws.Range("M22").Value = -51341.5
$ws.Range("H27").Value = 33869
$ws.Range("I27").Value = 51636.5
$ws.Range("K27").Value = 51636.5
$ws.Range("M27").Value = -51529.5
$ws.Range("H46").Value = 3338.8
$ws.Range("I46").Value = 1436.6666
$ws.Range("K46").Value = 1436.6666
$ws.Range("M46").Value = -1248.6666
$ws.Range("H55").Value = 3334091
$ws.Range("J55").Value = 1125.125
$ws.Range("L55").Value = 1125.125
$ws.Range("N55").Value = -1471.125
$ws.Range("H68").Value = 3656.0952
$ws.Range("I68").Value = 2946.0527
$ws.Range("K68").Value = 2946.0527
$ws.Range("M68").Value = -2197.0527
$ws.Range("H71").Value = 3656.0952
$ws.Range("I71").Value = 2946.0527
$ws.Range("K71").Value = 14730.2635
$ws.Range("M71").Value = -10986.2635
$ws.Range("H100").Value = 4385
$ws.Range("I100").Value = 1948.5
$ws.Range("K100").Value = 1948.5
$ws.Range("M100").Value = -1407.5
$ws.Range("H132").Value = 4568.7666
$ws.Range("I132").Value = 3250
$ws.Range("J132").Value = 6546.9165
$ws.Range("K132").Value = 9750
$ws.Range("L132").Value = 19640.7495
$ws.Range("M132").Value = -7220
$ws.Range("N132").Value = -24700.7495
$ws.Range("H136").Value = 7095.9585
$ws.Range("I136").Value = 3486.2
$ws.Range("K136").Value = 10458.6
$ws.Range("M136").Value = -7908.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3242.4443
$ws.Range("I81").Value = 3242.4443
$ws.Range("K81").Value = 6484.8886
$ws.Range("M81").Value = -5423.8886
$ws.Range("H84").Value = 3242.4443
$ws.Range("I84").Value = 3242.4443
$ws.Range("K84").Value = 32424.443
$ws.Range("M84").Value = -27120.443
$ws.Range("H100").Value = 797.2727
$ws.Range("J100").Value = 1129.3334
$ws.Range("L100").Value = 2258.6668
$ws.Range("N100").Value = -3340.6668
$ws.Range("H107").Value = 1963
$ws.Range("I107").Value = 2500.6667
$ws.Range("J107").Value = 350
$ws.Range("K107").Value = 7502.000100000001
$ws.Range("L107").Value = 1050
$ws.Range("M107").Value = -5582.000100000001
$ws.Range("N107").Value = -4890
$ws.Range("H132").Value = 6242.086
$ws.Range("I132").Value = 2882.3
$ws.Range("K132").Value = 8646.900000000001
$ws.Range("M132").Value = -6116.900000000001
